$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 8.426422666666666
$ws.Range("H2").Value = 25.279268
$ws.Range("I2").Value = 0.1151758588783328
$ws.Range("J2").Value = 0.1151758588783328
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 22.906497
$ws.Range("N2").Value = 68.719491
$ws.Range("O2").Value = 0.9446038650914245
$ws.Range("P2").Value = 0.9446038650914245
$ws.Range("Q2").Value = 193.019825534732
$ws.Range("R2").Value = 1737.178429812588
$ws.Range("S2").Value = 0.1087955614616976
$ws.Range("T2").Value = 0.1087955614616976

# Row 3
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 8.426422666666666
$ws.Range("H3").Value = 25.279268
$ws.Range("I3").Value = 0.1151758588783328
$ws.Range("J3").Value = 0.1151758588783328
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1329193333333333
$ws.Range("N3").Value = 0.3987579999999999
$ws.Range("O3").Value = 0.005481244732096839
$ws.Range("P3").Value = 0.005481244732096839
$ws.Range("Q3").Value = 1.120034483238222
$ws.Range("R3").Value = 10.080310349144
$ws.Range("S3").Value = 0.0006313070697415905
$ws.Range("T3").Value = 0.0006313070697415905

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 8.426422666666666
$ws.Range("H4").Value = 25.279268
$ws.Range("I4").Value = 0.1151758588783328
$ws.Range("J4").Value = 0.1151758588783328
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.210428333333333
$ws.Range("N4").Value = 3.631285
$ws.Range("O4").Value = 0.04991489017647865
$ws.Range("P4").Value = 0.04991489017647865
$ws.Range("Q4").Value = 10.19958074437555
$ws.Range("R4").Value = 91.79622669938
$ws.Range("S4").Value = 0.005748990346893585
$ws.Range("T4").Value = 0.005748990346893584

# Row 5
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 11.78712033333333
$ws.Range("H5").Value = 35.361361
$ws.Range("I5").Value = 0.1611112760180311
$ws.Range("J5").Value = 0.1611112760180311
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 22.906497
$ws.Range("N5").Value = 68.719491
$ws.Range("O5").Value = 0.9446038650914245
$ws.Range("P5").Value = 0.9446038650914245
$ws.Range("Q5").Value = 270.0016365541391
$ws.Range("R5").Value = 2430.014728987251
$ws.Range("S5").Value = 0.1521863340364435
$ws.Range("T5").Value = 0.1521863340364435

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 11.78712033333333
$ws.Range("H6").Value = 35.361361
$ws.Range("I6").Value = 0.1611112760180311
$ws.Range("J6").Value = 0.1611112760180311
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1329193333333333
$ws.Range("N6").Value = 0.3987579999999999
$ws.Range("O6").Value = 0.005481244732096839
$ws.Range("P6").Value = 0.005481244732096839
$ws.Range("Q6").Value = 1.566736176626444
$ws.Range("R6").Value = 14.100625589638
$ws.Range("S6").Value = 0.0008830903329552329
$ws.Range("T6").Value = 0.0008830903329552329

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 11.78712033333333
$ws.Range("H7").Value = 35.361361
$ws.Range("I7").Value = 0.1611112760180311
$ws.Range("J7").Value = 0.1611112760180311
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.210428333333333
$ws.Range("N7").Value = 3.631285
$ws.Range("O7").Value = 0.04991489017647865
$ws.Range("P7").Value = 0.04991489017647865
$ws.Range("Q7").Value = 14.26746441987611
$ws.Range("R7").Value = 128.407179778885
$ws.Range("S7").Value = 0.008041851648632362
$ws.Range("T7").Value = 0.008041851648632362

# Row 8 (new)
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Vcam1"
$ws.Range("C8").Value = "Itga4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 52.94781866666667
$ws.Range("H8").Value = 158.843456
$ws.Range("I8").Value = 0.7237128651036362
$ws.Range("J8").Value = 0.7237128651036362
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 22.906497
$ws.Range("N8").Value = 68.719491
$ws.Range("O8").Value = 0.9446038650914245
$ws.Range("P8").Value = 0.9446038650914245
$ws.Range("Q8").Value = 1212.849049444544
$ws.Range("R8").Value = 10915.6414450009
$ws.Range("S8").Value = 0.6836219695932835
$ws.Range("T8").Value = 0.6836219695932835

# Row 9 (new)
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Vcam1"
$ws.Range("C9").Value = "Itga4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 52.94781866666667
$ws.Range("H9").Value = 158.843456
$ws.Range("I9").Value = 0.7237128651036362
$ws.Range("J9").Value = 0.7237128651036362
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1329193333333333
$ws.Range("N9").Value = 0.3987579999999999
$ws.Range("O9").Value = 0.005481244732096839
$ws.Range("P9").Value = 0.005481244732096839
$ws.Range("Q9").Value = 7.037788758627554
$ws.Range("R9").Value = 63.34009882764799
$ws.Range("S9").Value = 0.003966847329400017
$ws.Range("T9").Value = 0.003966847329400017

# Row 10 (new)
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Vcam1"
$ws.Range("C10").Value = "Itga4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 52.94781866666667
$ws.Range("H10").Value = 158.843456
$ws.Range("I10").Value = 0.7237128651036362
$ws.Range("J10").Value = 0.7237128651036362
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.210428333333333
$ws.Range("N10").Value = 3.631285
$ws.Range("O10").Value = 0.04991489017647865
$ws.Range("P10").Value = 0.04991489017647865
$ws.Range("Q10").Value = 64.0895399023289
$ws.Range("R10").Value = 576.80585912096
$ws.Range("S10").Value = 0.03612404818095271
$ws.Range("T10").Value = 0.03612404818095271
